$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:B25").ClearContents()
$ws.Range("A13").Select()
